# Renumber the "Listing" index in column A of the Summary sheet.
# The author's fix replaces the old (non-contiguous) listing numbers in
# column A with a simple contiguous sequence starting at 0 (row 2 -> 0,
# row 3 -> 1, row 4 -> 2, ... row 124 -> 122).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 124
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
